# Elixir - Day 3: insert a new "Homework" slide right before the
# "Q & A" slide. The new slide becomes slide 20 (id 347), pushing
# "Q & A" (id 287) to 21 and "Thank you" (id 299) to 22.

$p = $ppt.ActivePresentation

# Locate the "Q & A" slide so the new slide can be inserted right before
# it, regardless of its current absolute index.
$qaIndex = 0
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    if ($s.Shapes.Item(1).TextFrame.TextRange.Text -eq "Q & A") {
        $qaIndex = $i
        break
    }
}

# Borrow the "Title and Content" layout from an existing slide (slide 7,
# "Creating a new project") that already uses it, since the
# CustomLayouts collection can't reliably be indexed directly.
$titleContentLayout = $p.Slides.Item(7).CustomLayout

$newSlide = $p.Slides.AddSlide($qaIndex, $titleContentLayout)

# --- Title placeholder ---------------------------------------------------
$title = $newSlide.Shapes.Item(1)
$title.Name = "Title 3"
$title.TextFrame.TextRange.Text = "Homework"

# --- Body placeholder ------------------------------------------------------
$body = $newSlide.Shapes.Item(2)
$body.Name = "Text Placeholder 4"
$bodyRange = $body.TextFrame.TextRange

# Build up the paragraphs incrementally (rather than one bulk multi-line
# assignment) so each run keeps its "en-US" language tag.
$bodyRange.Text = "https://github.com/georgiyolovski/elixir-workshop/tree/main/day3/homework"
$bodyRange.Text = $bodyRange.Text + "`r"
$bodyRange.Text = $bodyRange.Text + "`rThe following module references might be helpful:"
$bodyRange.Text = $bodyRange.Text + "`rhttps://hexdocs.pm/phoenix/Mix.Tasks.Phx.New.html"
$bodyRange.Text = $bodyRange.Text + "`rhttps://hexdocs.pm/phoenix/Mix.Tasks.Phx.Gen.Json.html"

# Paragraph 1: link to the homework repo
$para1 = $bodyRange.Paragraphs(1)
$para1.ActionSettings(1).Hyperlink.Address = "https://github.com/georgiyolovski/elixir-workshop/tree/main/day3/homework"

# Paragraph 2 is left blank.
# Paragraph 3: "The following module references might be helpful:"

# Paragraph 4 & 5: indented references, each with its own hyperlink.
$para4 = $bodyRange.Paragraphs(4)
$para4.IndentLevel = 2
$para4.ActionSettings(1).Hyperlink.Address = "https://hexdocs.pm/phoenix/Mix.Tasks.Phx.New.html"

$para5 = $bodyRange.Paragraphs(5)
$para5.IndentLevel = 2
$para5.ActionSettings(1).Hyperlink.Address = "https://hexdocs.pm/phoenix/Mix.Tasks.Phx.Gen.Json.html"
